$wb = $excel.ActiveWorkbook

# 1. Add new workbook-level defined names
$wb.Names.Add("Attributes.Enabled", "='CAN1'!`$F`$2")
$wb.Names.Add("Attributes.ReadRate", "='CAN1'!`$E`$2")

# 2. Rename the F1 header on CAN1 from "Attributes.Boolean" to "Attributes.Enabled"
$canSheet = $wb.Worksheets.Item("CAN1")
$canSheet.Range("F1").Value = "Attributes.Enabled"

# 3. Update selection on CAN1 sheet to F2
[void]$canSheet.Range("F2").Select()

# 4. Activate CAN1 sheet (becomes the selected/active tab)
$canSheet.Activate()
